$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hồ sơ nhân sự")

# Column D holds properties.last_edited_time for each Notion page row.
# Rows 2-18 previously shared the value 2024-08-03T03:18:00.000Z, which is
# now refreshed to 2024-08-03T03:29:00.000Z.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T03:29:00.000Z"
}

# Rows 19-20 (Trương Lâm Khanh, Đỗ Thị Huyền Trân) get the newly edited
# last_edited_time value.
$ws.Cells.Item(19, 4).Value = "2024-08-03T03:28:00.000Z"
$ws.Cells.Item(20, 4).Value = "2024-08-03T03:28:00.000Z"

# Rows 21-26 previously shared the value 2024-08-03T03:17:00.000Z, which is
# now refreshed to 2024-08-03T03:28:00.000Z.
for ($r = 21; $r -le 26; $r++) {
    $ws.Cells.Item($r, 4).Value = "2024-08-03T03:28:00.000Z"
}
